$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 entirely (shifts nothing up since they're the last rows,
# but removes them from the sheet dimension / data)
$ws.Range("A4:H5").EntireRow.Delete()

# Update row 2 values
$ws.Range("B2").Value = 75000
$ws.Range("D2").Value = 0.3803303190246927
$ws.Range("E2").Value = 3.602579562328429
$ws.Range("F2").Value = 0.7518
$ws.Range("H2").Value = 3.748694078222786

# Update row 3 values
$ws.Range("B3").Value = 75000
$ws.Range("D3").Value = 0.6900868030522065
$ws.Range("E3").Value = 3.363668530993691
$ws.Range("F3").Value = 1.3849
$ws.Range("H3").Value = 3.748694078222786
